$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before DB (existing DB -> DC, DC -> DD, etc.)
$ws.Range("DB:DB").Insert()

# Row 1 header: new snapshot timestamp for the freshly inserted column
$ws.Range("DB1").Value = "2026-02-01 13:36:07"

# Data rows: the new column DB is a duplicate snapshot of the previous
# last-price column (DA), same as every other timestamp column that was
# appended before it.
$ws.Range("DA2:DA206").Copy($ws.Range("DB2:DB206"))

Write-Output "done"
